$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25, pushing existing rows (and their
# formatting) down by one. This matches the diff: every row from the
# old 25..153 becomes 26..154, and a brand new row of data is placed
# at row 25.
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with this week's data.
$ws.Range("A25").Value = 4
$ws.Range("B25").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C25").Value = "Los Lagos"
$ws.Range("D25").Value = 44547
$ws.Range("E25").Value = 10
$ws.Range("F25").Value = 100112039
$ws.Range("G25").Value = "Ciboulette"
$ws.Range("H25").Value = "Sin especificar"
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = 2500
$ws.Range("N25").Value = "$/docena de atados"
$ws.Range("O25").Value = "Región Metropolitana"
$ws.Range("P25").Value = 833
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = "Hortaliza"
